# "change in detail page"
#
# This script:
#  1. Rewrites the "Please login on contact me (in offline mode) " paragraph:
#     splits it into several runs (with a duplicated/garbled "in" and new
#     "offer detail" text), colours it green (00B050) and adds the
#     proof-reading marks Word inserts around the new grammar break.
#  2. Inserts a new, empty (bold-only) paragraph right after it.
#  3. Colours a handful of other paragraphs on the same "detail page" topic
#     green (00B050), run by run, leaving their text/proofErr markers as-is.
#  4. Strips the stray <w:lastRenderedPageBreak/> from the
#     "12. The calendar ..." run.

$d = $word.ActiveDocument

# wdColor value for RGB 00B050 (packed as 0x00BBGGRR -> 0x0050B000)
$green = 5287936

# ---------------------------------------------------------------------
# 1. "Please login on contact me (in offline mode) " -> split into runs,
#    reworded, coloured green.
#    (Paragraph.Range.Text always carries a trailing paragraph-mark
#    character, so match with -like / a trailing "*" rather than -eq.)
# ---------------------------------------------------------------------
$pleaseXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:b/><w:color w:val="00B050"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:color w:val="00B050"/></w:rPr><w:t>Please login on contact me (in offline mode</w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:b/><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve"> in</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:b/><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve"> offer detail </w:t></w:r>' +
  '</w:p>'

$pleasePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Please login on contact me (in offline mode) *") {
        $pleasePara = $p
        break
    }
}
$pleasePara.Range.InsertXML($pleaseXml)

# ---------------------------------------------------------------------
# 2. Insert a brand-new empty paragraph (bold only, no colour) right
#    after it.
# ---------------------------------------------------------------------
$pleasePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Please login on contact me*offer detail *") {
        $pleasePara = $p
        break
    }
}
$insertPos = $pleasePara.Range.End - 1
$insertRng = $d.Range($insertPos, $insertPos)
$insertRng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>')

# ---------------------------------------------------------------------
# 3. Colour a batch of whole paragraphs green (00B050), run by run -
#    this preserves existing proofErr markers / tabs / text untouched
#    and only adds <w:color w:val="00B050"/> to each run's (and the
#    paragraph mark's) rPr.
# ---------------------------------------------------------------------
$targets = @(
    "User store button (go to offer owner store) rather owner",
    "No video url mandatory ",
    "And seo title,seo_description and keywords (whole section shuldnt need)",
    "Last check box should be mandatory",
    "For  sound if not found no element should be shown)"
)

foreach ($targetText in $targets) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "$targetText*") {
            $p.Range.Font.Color = $green
            break
        }
    }
}

# ---------------------------------------------------------------------
# 4. Remove the stray <w:lastRenderedPageBreak/> in front of
#    "12. The calendar under Price Calculation should be OUT".
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "12. The calendar under Price Calculation should be OUT*") {
        $calendarXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
          '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="00B050"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr>' +
          '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="00B050"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>' +
          '<w:t>12. The calendar under Price Calculation should be OUT</w:t></w:r></w:p>'
        $p.Range.InsertXML($calendarXml)
        break
    }
}

Write-Host "done"
